$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 595 -- this shifts rows 595:631 down to 596:632
# (and their formatting/styles move with them, like a real Excel "Insert" does).
$ws.Rows.Item(595).Insert()

# Fill in the new row 595 with this week's data. Columns A,B,C,E,F,G,H,I,J,K,L,Q,T
# are identical to the surrounding "Femacal de La Calera - Mango" rows; only the
# date + volume/price/origin columns for this record are new.
$ws.Cells.Item(595, 1).Value = 3
$ws.Cells.Item(595, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(595, 3).Value = "Coquimbo"
$ws.Cells.Item(595, 4).Value = 45041
$ws.Cells.Item(595, 5).Value = 5
$ws.Cells.Item(595, 6).Value = "Fruta"
$ws.Cells.Item(595, 7).Value = 100108
$ws.Cells.Item(595, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(595, 9).Value = 100108002
$ws.Cells.Item(595, 10).Value = "Mango"
$ws.Cells.Item(595, 11).Value = "Sin especificar"
$ws.Cells.Item(595, 12).Value = "Primera"
$ws.Cells.Item(595, 13).Value = 228
$ws.Cells.Item(595, 14).Value = 7000
$ws.Cells.Item(595, 15).Value = 7000
$ws.Cells.Item(595, 16).Value = 7000
$ws.Cells.Item(595, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(595, 18).Value = "Perú"
$ws.Cells.Item(595, 19).Value = 1750
$ws.Cells.Item(595, 20).Value = 4
